$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.740.26"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.871.62"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'336.49"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").Value = "'0.4686"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").Value = "'0.3928"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.08015"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "'45.19"
$ws.Range("E10").Value = "  -5.08%  "
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").Value = "'21.89"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "1.890.50"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "'6.004"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "'7.277"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "'88.89"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "'0.06760"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "'1.012"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").Value = "27.761.90"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'5.495"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "'10.95"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Value = "'2.312"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").Value = "2.100.68"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'159.80"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "'19.83"
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").Value = "'2.148"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").Value = "'121.94"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "'0.9829"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'0.09525"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").Value = "'3.642"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").Value = "'5.346"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'1.341"
$ws.Range("E36").Value = "  -7.58%  "
$ws.Range("D37").Value = "'0.06072"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'0.02243"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").Value = "'8.304"
$ws.Range("E40").Value = "  +3.02%  "
$ws.Range("D41").Value = "'1.011"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "'0.5998"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'0.1894"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "'1.245"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").Value = "'0.5674"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "'12.27"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'0.06770"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("D50").Value = "'112.40"
$ws.Range("D51").Value = "'3.022"
$ws.Range("E51").Value = "  -11.03%  "
